$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.929.20"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.48%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.554.38"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.30%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.41%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.85"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.28%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.574"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.42%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.09%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.540"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.26%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.71"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.92%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0812"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.11%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.51"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.92%  "

$ws.Range("E13").Value = "  -4.90%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.941.24"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.19%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.536.96"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.72%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.09"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.29%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.850"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.05%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.928.15"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.24%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.86"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.32%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.56"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.08%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0961"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.55%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.72"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.09%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "253.56"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.60%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.96"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.93%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.06"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.12%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.79"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.02%  "

$ws.Range("E27").Value = "  +0.31%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.41"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.48%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "40.47"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.51%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.27"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.78%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.86"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.59%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "156.54"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.44%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.51"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.77%  "

$ws.Range("E34").Value = "  +3.35%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0801"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.86%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.11"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.14%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.31"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.37%  "

$ws.Range("E38").Value = "  +1.44%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.45"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.46%  "

$ws.Range("E40").Value = "  +0.28%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.02"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -7.24%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.82"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.97%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0304"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.90%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.01%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.26"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.32%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.988.66"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.25%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "84.54"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.53%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.04"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.87%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.790.93"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.01%  "

$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "104.86"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.22%  "

$ws.Range("B51").Value = "ordi"
$ws.Range("C51").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "74.43"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.15%  "
